$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    $cell.NumberFormat = "@"
    $cell.Value = "2013-06-25"
}
